$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values (TB=B, d2S=C, K=D, IP=E, sum=G) for rows 2-13, regenerated
# s_val data after filtering save games.
$data = @{
    2  = @(0.6606524410359556, 0.306821227259698,  0.1494219747398047, 0.4942365360607697, 1.611132179096228)
    3  = @(0.1190320826869504, 1.655778082260271,  3.537761648806719,  0.4942365360607697, 5.806808349814711)
    4  = @(1.455362044514542,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    5  = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    6  = @(1.455362044514542,  1.655778082260271,  0.1494219747398047, 0.4942365360607697, 3.754798637575387)
    7  = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    8  = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    9  = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    10 = @(3.286832544864788,  1.655778082260271,  0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    11 = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    12 = @(0.6606524410359556, 1.655778082260271,  0.1494219747398047, 10.19245300693656,  12.65830550497259)
    13 = @(3.286832544864788,  1.655778082260271,  3.537761648806719,  10.19245300693656,  18.67282528286833)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Cells.Item($r, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($r, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($r, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($r, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($r, 7).Value = $vals[4]   # G - sum
}
